$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.305.07"
$ws.Range("E2").Value = "  -2.74%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.301.46"
$ws.Range("E3").Value = "  -3.38%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "557.48"
$ws.Range("E5").Value = "  -3.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.98"
$ws.Range("E6").Value = "  -4.59%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.301.75"
$ws.Range("E8").Value = "  -3.39%  "
$ws.Range("E9").Value = "  -3.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.84"
$ws.Range("E10").Value = "  -2.67%  "
$ws.Range("E11").Value = "  -3.70%  "
$ws.Range("E12").Value = "  -1.55%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.870.80"
$ws.Range("E13").Value = "  -3.29%  "
$ws.Range("E14").Value = "  +0.15%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.81"
$ws.Range("E15").Value = "  -5.55%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.307.48"
$ws.Range("E16").Value = "  -3.31%  "
$ws.Range("E17").Value = "  -3.00%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "60.339.87"
$ws.Range("E18").Value = "  -2.70%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.10"
$ws.Range("E19").Value = "  -4.29%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.10"
$ws.Range("E20").Value = "  -2.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.64"
$ws.Range("E21").Value = "  -2.54%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "374.88"
$ws.Range("E22").Value = "  -1.39%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "74.47"
$ws.Range("E23").Value = "  -0.93%  "
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.536"
$ws.Range("E25").Value = "  -4.95%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.445.57"
$ws.Range("E26").Value = "  -3.10%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000102"
$ws.Range("E27").Value = "  -7.86%  "
$ws.Range("E28").Value = "  -4.38%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  -0.31%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.17"
$ws.Range("E30").Value = "  -6.02%  "
$ws.Range("E31").Value = "  -0.04%  "
$ws.Range("E32").Value = "  -3.63%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.57"
$ws.Range("E33").Value = "  -4.08%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "22.60"
$ws.Range("E34").Value = "  -2.03%  "
$ws.Range("E35").Value = "  -6.72%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.11"
$ws.Range("E36").Value = "  -6.30%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "166.69"
$ws.Range("E37").Value = "  -1.14%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.52"
$ws.Range("E38").Value = "  -3.74%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.69"
$ws.Range("E39").Value = "  -2.48%  "
$ws.Range("B40").Value = "EnergySwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "26.87"
$ws.Range("E40").Value = "  -13.37%  "
$ws.Range("B41").Value = "RenzoRestakedETH"
$ws.Range("C41").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.336.13"
$ws.Range("E41").Value = "  -3.28%  "
$ws.Range("E42").Value = "  -6.48%  "
$ws.Range("E43").Value = "  -1.90%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.751"
$ws.Range("E44").Value = "  -3.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.14"
$ws.Range("E45").Value = "  -4.98%  "
$ws.Range("E46").Value = "  -5.48%  "
$ws.Range("E47").Value = "  -4.06%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.353.17"
$ws.Range("E48").Value = "  -7.51%  "
$ws.Range("E49").Value = "  +0.05%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.41"
$ws.Range("E50").Value = "  -6.56%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.31"
$ws.Range("E51").Value = "  -4.76%  "
